$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update row 2 and row 7 values in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5706
$ws1.Range("F7").Value = 69

# Sheet "全部类型" (all types) - same underlying rows updated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5706
$ws4.Range("F7").Value = 69
